$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (not be coerced to a
# number). Type it with a leading apostrophe in a scratch cell, then
# copy/paste-values into B3 so the original cell style (s=8) is preserved
# instead of Excel minting a new "quotePrefix" style directly on B3.
$ws.Range("Z1").Value = "'2570314725427075"
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()

$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 22.02.2025"

# Row 6
$ws.Range("B6").Value = "24.02."
$ws.Range("C6").Value = "25.02."
$ws.Range("D6").Value = "MCDONALDS Aachen"
$ws.Range("E6").Value = "13,92-"

# Row 7
$ws.Range("B7").Value = "28.02."
$ws.Range("C7").Value = "01.03."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "24,85-"

# Row 8
$ws.Range("B8").Value = "01.03."
$ws.Range("C8").Value = "02.03."
$ws.Range("D8").Value = "ZALANDO MKTPLC EU ATKHSW"
$ws.Range("E8").Value = "72,11-"

# Row 9: transaction removed, row now blank (E9 keeps a centered style)
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E9").VerticalAlignment = -4108    # xlCenter
$ws.Range("E9").WrapText = $true

# Row 10: transaction removed, row now blank (E10 keeps a right-aligned style)
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152  # xlRight
$ws.Range("E10").VerticalAlignment = -4108    # xlCenter
$ws.Range("E10").WrapText = $true

# Row 12: closing balance date/amount
$ws.Range("D12").Value = "KONTOSTAND AM 04.03.2025"
$ws.Range("E12").Value = "110,88-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 12.03.2025"
